$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 - "2021年" data row appended below the existing 2011-2020 rows.
$ws.Cells.Item(12, 1).Value = "2021年"
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 6).Value = 204
$ws.Cells.Item(12, 7).Value = 29
$ws.Cells.Item(12, 8).Value = 41
$ws.Cells.Item(12, 9).Value = 26
$ws.Cells.Item(12, 14).Value = 172
$ws.Cells.Item(12, 15).Value = 103
$ws.Cells.Item(12, 16).Value = 12
$ws.Cells.Item(12, 17).Value = 17
$ws.Cells.Item(12, 18).Value = 22
$ws.Cells.Item(12, 20).Value = 12
$ws.Cells.Item(12, 24).Value = 40
$ws.Cells.Item(12, 25).Value = 107
$ws.Cells.Item(12, 28).Value = 11747
$ws.Cells.Item(12, 30).Value = 10
$ws.Cells.Item(12, 31).Value = 18
$ws.Cells.Item(12, 34).Value = 12
$ws.Cells.Item(12, 36).Value = 162
$ws.Cells.Item(12, 38).Value = 76
$ws.Cells.Item(12, 39).Value = 351
$ws.Cells.Item(12, 41).Value = 4642
$ws.Cells.Item(12, 42).Value = 15
$ws.Cells.Item(12, 45).Value = 22
$ws.Cells.Item(12, 48).Value = 5651

# A12 ("2021年") should carry the same formatting as the other year labels
# in column A (bold, centered, bordered) - copy the style from A11.
$ws.Cells.Item(11, 1).Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Row 12 written"
